$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Advance the Instance test case name from Automation2 to Automation3
$ws.Range("D2").Value = "Automation3"

# Move the selection to E2 (carry forward hourly run)
$ws.Range("E2").Select()
